# Additional transportation calibration edits to adjust psgr LDV lifetime
# and moving TTS s-curve for BEVs to reach 1 by 2035

$wb = $excel.ActiveWorkbook

# Work on the passenger LDV sheet
$ws = $wb.Worksheets.Item("SoCDTtiNTY-psgr")

# Row 2 = LDVs. Set the BEV share-that-is-new s-curve input (B2) to 1/20,
# and fill the shared formula across C2:H2 (natural gas, gasoline, diesel,
# plugin hybrid, LPG, hydrogen vehicles) so it reaches 1 by 2035.
$ws.Range("B2").Formula = "=1/20"
$ws.Range("C2:H2").Formula = "=1/20"

# Make this sheet the active/selected one, matching the new selection
$ws.Activate()
$ws.Range("B2:H2").Select()
